$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values in column F (dSF) to match repulled/recalculated data.
# Map of row -> new value for column F
$updates = @{
    2  = -7
    3  = -6
    9  = -2
    10 = -2
    15 = 0
    19 = 0
    20 = -6
    22 = 4
    27 = -1
    29 = -3
    37 = -9
    38 = -3
    41 = 12
    42 = 1
    53 = -6
    62 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
